$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 7: set "Reported by" (column E) to Adrian
$ws.Range("E7").Value = "Adrian"

# Row 8: new known issue entry
$ws.Range("A8").Value = Get-Date -Year 2015 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B8").Value = "medium"
$ws.Range("C8").Value = "task controller and Vchans"
$ws.Range("D8").Value = "When removing a task controller from the task tree, its source Vchans are not sending data anymore to SinkVChans belonging to the task tree, and the task tree execution will time out. To correct this, Source Vchans should be also associated (and registered) with task controllers so that they are disconnected from their sinks if the task controller is removed from the task tree"
$ws.Range("E8").Value = "Adrian"

$ws.Rows.Item(8).RowHeight = 60

# Update selection to A9
$ws.Range("A9").Select()
